$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    2  = 463.3
    3  = 472.3
    4  = 477
    5  = 463.1
    6  = 458
    7  = 432.7
    8  = 428.9
    9  = 427
    10 = 435.7
    11 = 435.2
    12 = 415.2
    13 = 406.8
    14 = 406
    15 = 415.3
    16 = 423.2
    17 = 422.7
    18 = 433
    19 = 431.7
    20 = 427.9
    21 = 431.7
}

foreach ($row in $values.Keys) {
    $ws.Range("C$row").Value = $values[$row]
}
